$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 11715.3
$ws.Range("J17").Value = 12794.777
$ws.Range("L17").Value = 38384.331
$ws.Range("N17").Value = -38720.331

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 241.625
$ws.Range("I28").Value = 204.46666
$ws.Range("K28").Value = 204.46666
$ws.Range("M28").Value = 280.53334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 10403.565
$ws.Range("I62").Value = 11595.091
$ws.Range("K62").Value = 11595.091
$ws.Range("M62").Value = -10971.091

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4987.8
$ws.Range("J64").Value = 4987.8
$ws.Range("L64").Value = 4987.8
$ws.Range("N64").Value = -5483.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 10403.565
$ws.Range("I65").Value = 11595.091
$ws.Range("K65").Value = 57975.455
$ws.Range("M65").Value = -54855.455

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4987.8
$ws.Range("J67").Value = 4987.8
$ws.Range("L67").Value = 4987.8
$ws.Range("N67").Value = -6703.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2293.5
$ws.Range("I88").Value = 2890
$ws.Range("J88").Value = 2094.6667
$ws.Range("K88").Value = 2890
$ws.Range("L88").Value = 2094.6667
$ws.Range("M88").Value = -2484
$ws.Range("N88").Value = -2906.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2293.5
$ws.Range("I91").Value = 2890
$ws.Range("J91").Value = 2094.6667
$ws.Range("K91").Value = 2890
$ws.Range("L91").Value = 2094.6667
$ws.Range("M91").Value = -1486
$ws.Range("N91").Value = -4902.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1750.4445
$ws.Range("J111").Value = 482
$ws.Range("L111").Value = 1446
$ws.Range("N111").Value = -7580

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 109998
$ws.Range("J128").Value = 109998
$ws.Range("L128").Value = 109998
$ws.Range("N128").Value = -119958

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19646.555
$ws.Range("I32").Value = 21253.574
$ws.Range("K32").Value = 21253.574
$ws.Range("M32").Value = -20966.574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3758
$ws.Range("I45").Value = 2998.6191
$ws.Range("K45").Value = 2998.6191
$ws.Range("M45").Value = -2621.6191

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2082.5
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2251.8462
$ws.Range("I88").Value = 1699.5
$ws.Range("J88").Value = 2352.2727
$ws.Range("K88").Value = 1699.5
$ws.Range("L88").Value = 2352.2727
$ws.Range("M88").Value = -1293.5
$ws.Range("N88").Value = -3164.2727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2251.8462
$ws.Range("I91").Value = 1699.5
$ws.Range("J91").Value = 2352.2727
$ws.Range("K91").Value = 1699.5
$ws.Range("L91").Value = 2352.2727
$ws.Range("M91").Value = -295.5
$ws.Range("N91").Value = -5160.2727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2276.25
$ws.Range("I102").Value = 1983.4348
$ws.Range("K102").Value = 1983.4348
$ws.Range("M102").Value = -361.4348

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 45920.523
$ws.Range("I132").Value = 49981.145
$ws.Range("J132").Value = 3284
$ws.Range("K132").Value = 149943.435
$ws.Range("L132").Value = 9852
$ws.Range("M132").Value = -147413.435
$ws.Range("N132").Value = -14912

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2082.5
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 95599.5
$ws.Range("J57").Value = 95599.5
$ws.Range("L57").Value = 95599.5
$ws.Range("N57").Value = -97039.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3459.111
$ws.Range("I105").Value = 3484.4
$ws.Range("J105").Value = 3332.6667
$ws.Range("K105").Value = 3484.4
$ws.Range("L105").Value = 3332.6667
$ws.Range("M105").Value = -1737.4
$ws.Range("N105").Value = -6826.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1072.1666
$ws.Range("I107").Value = 715.0909
$ws.Range("K107").Value = 715.0909
$ws.Range("M107").Value = 1204.9091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2359.1428
$ws.Range("I134").Value = 1819
$ws.Range("J134").Value = 5600
$ws.Range("K134").Value = 5457
$ws.Range("L134").Value = 16800
$ws.Range("M134").Value = -2922
$ws.Range("N134").Value = -21870

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 95599.5
$ws.Range("J136").Value = 95599.5
$ws.Range("L136").Value = 95599.5
$ws.Range("N136").Value = -105799.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 81200
$ws.Range("J141").Value = 81200
$ws.Range("L141").Value = 81200
$ws.Range("N141").Value = -91560

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1784.0769
$ws.Range("J22").Value = 3566.1667
$ws.Range("L22").Value = 3566.1667
$ws.Range("N22").Value = -4266.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1585.5
$ws.Range("I31").Value = 1566.8422
$ws.Range("K31").Value = 1566.8422
$ws.Range("M31").Value = -1271.8422

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1585.5
$ws.Range("I34").Value = 1566.8422
$ws.Range("K34").Value = 1566.8422
$ws.Range("M34").Value = -1364.8422

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3641.3
$ws.Range("I99").Value = 4120.6665
$ws.Range("K99").Value = 4120.6665
$ws.Range("M99").Value = -2622.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1695
$ws.Range("I105").Value = 1656.875
$ws.Range("K105").Value = 1656.875
$ws.Range("M105").Value = 90.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3641.3
$ws.Range("I126").Value = 4120.6665
$ws.Range("K126").Value = 12361.9995
$ws.Range("M126").Value = -9891.999500000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 51081.906
$ws.Range("I134").Value = 69047.13
$ws.Range("K134").Value = 207141.39
$ws.Range("M134").Value = -204606.39

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 387.9
$ws.Range("I114").Value = 387.9
$ws.Range("K114").Value = 1163.7
$ws.Range("M114").Value = 2090.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6268344.5
$ws.Range("I131").Value = 14272.875
$ws.Range("J131").Value = 12522416
$ws.Range("K131").Value = 42818.625
$ws.Range("L131").Value = 37567248
$ws.Range("M131").Value = -37778.625
$ws.Range("N131").Value = -37577328

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 592.087
$ws.Range("I97").Value = 496.5
$ws.Range("K97").Value = 496.5
$ws.Range("M97").Value = -0.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7884.375
$ws.Range("I126").Value = 4695
$ws.Range("J126").Value = 8947.5
$ws.Range("K126").Value = 14085
$ws.Range("L126").Value = 26842.5
$ws.Range("M126").Value = -11615
$ws.Range("N126").Value = -31782.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11784.059
$ws.Range("J7").Value = 5386.8887
$ws.Range("L7").Value = 5386.8887
$ws.Range("N7").Value = -5610.8887

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 16387.408
$ws.Range("I46").Value = 23612.428
$ws.Range("J46").Value = 3743.625
$ws.Range("K46").Value = 23612.428
$ws.Range("L46").Value = 3743.625
$ws.Range("M46").Value = -23424.428
$ws.Range("N46").Value = -4119.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2108
$ws.Range("I82").Value = 771
$ws.Range("J82").Value = 3445
$ws.Range("K82").Value = 771
$ws.Range("L82").Value = 3445
$ws.Range("M82").Value = -410
$ws.Range("N82").Value = -4167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2108
$ws.Range("I85").Value = 771
$ws.Range("J85").Value = 3445
$ws.Range("K85").Value = 771
$ws.Range("L85").Value = 3445
$ws.Range("M85").Value = 477
$ws.Range("N85").Value = -5941

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4061.926
$ws.Range("I122").Value = 3341.5715
$ws.Range("K122").Value = 10024.7145
$ws.Range("M122").Value = -7574.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 11784.059
$ws.Range("J126").Value = 5386.8887
$ws.Range("L126").Value = 16160.6661
$ws.Range("N126").Value = -21100.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3166.6667
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 4250
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 4250
$ws.Range("M2").Value = -888
$ws.Range("N2").Value = -4474

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 24038
$ws.Range("I49").Value = 24038
$ws.Range("K49").Value = 24038
$ws.Range("M49").Value = -23808

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 187730.4
$ws.Range("I126").Value = 207144.33
$ws.Range("J126").Value = 13005
$ws.Range("K126").Value = 621432.99
$ws.Range("L126").Value = 39015
$ws.Range("M126").Value = -618962.99
$ws.Range("N126").Value = -43955

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 70907.21000000001
$ws.Range("I132").Value = 79016.08
$ws.Range("J132").Value = 3333.3333
$ws.Range("K132").Value = 237048.24
$ws.Range("L132").Value = 9999.999899999999
$ws.Range("M132").Value = -234518.24
$ws.Range("N132").Value = -15059.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1885.3334
$ws.Range("I136").Value = 1662.2941
$ws.Range("J136").Value = 2833.25
$ws.Range("K136").Value = 4986.8823
$ws.Range("L136").Value = 8499.75
$ws.Range("M136").Value = -2436.8823
$ws.Range("N136").Value = -13599.75
